# Punto 7 Programa 1
# - Rename header H1 from "Pregunta7" to "Mayor dif. PTS"
# - Fill in the "Mayor dif. PTS" column (H2:H14) with the per-team values
# - Update the "Mas dif Puntos:" summary row (row 23) to point at Valencia (the
#   max value) and change its separator text from " / " to " - "
# - Change the " / " separators to " - " in the other summary rows (17-22)
# - Add a new block (rows 25-39): "Mayor diferencia en temporada" table listing,
#   for every team, its max PTS difference and in which season-transition it
#   occurred

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 is a blank spacer row in the original sheet (empty string-typed
# cells with no value). Re-asserting empty strings here keeps it genuinely
# blank across the load/save round-trip.
$ws.Range("A15").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""

# --- Header ---------------------------------------------------------------
$ws.Range("H1").Value = "Mayor dif. PTS"

# --- "Mayor dif. PTS" column (H2:H14) --------------------------------------
$difPts = @{
    2  = 25
    3  = 20
    4  = 13
    5  = 7
    6  = 5
    7  = 6
    8  = 11
    9  = 12
    10 = 11
    11 = 15
    12 = 19
    13 = 24
    14 = 33
}
foreach ($row in $difPts.Keys) {
    $ws.Cells.Item($row, 8).Value = $difPts[$row]
}

# --- Summary rows: " / " -> " - " ------------------------------------------
$ws.Range("B17").Value = " - Barcelona"
$ws.Range("B18").Value = " - Real Sociedad"
$ws.Range("B19").Value = " - Rayo Vallecano"
$ws.Range("B20").Value = " - Real Madrid"
$ws.Range("B21").Value = " - Atlético Madrid"
$ws.Range("B22").Value = " - Barcelona - Real Madrid"

# --- "Mas dif Puntos:" row now reports the team with the biggest swing -----
$ws.Range("B23").Value = " - Valencia"
$ws.Range("C23").Value = 33

# --- New block: "Mayor diferencia en temporada" -----------------------------
# Row 24 is a blank spacer row, mirroring row 15's layout, before the new block.
$ws.Range("A24").Value = ""
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = ""

$ws.Range("A25").Value = "Mayor diferencia en temporada"

$ws.Range("A26").Value = "Equipo"
$ws.Range("B26").Value = "Dif"
$ws.Range("C26").Value = "Temporada"

$s1213_1314 = " 2012 - 2013 a 2013 - 2014/ "
$s1112_1213 = " 2011 - 2012 a 2012 - 2013/ "
$s1314_1415 = " 2013 - 2014 a 2014 - 2015/ "
$s1415_1516 = " 2014 - 2015 a 2015 - 2016/ "

$teamRows = @(
    @{ Row = 27; Team = "Athletic Bilbao";  Dif = 25; Seasons = @($s1213_1314) }
    @{ Row = 28; Team = "Atlético Madrid";  Dif = 20; Seasons = @($s1112_1213) }
    @{ Row = 29; Team = "Barcelona";        Dif = 13; Seasons = @($s1213_1314) }
    @{ Row = 30; Team = "Espanyol";         Dif = 7;  Seasons = @($s1314_1415) }
    @{ Row = 31; Team = "Getafe";           Dif = 5;  Seasons = @($s1213_1314, $s1314_1415) }
    @{ Row = 32; Team = "Granada";          Dif = 6;  Seasons = @($s1314_1415) }
    @{ Row = 33; Team = "Levante";          Dif = 11; Seasons = @($s1314_1415) }
    @{ Row = 34; Team = "Málaga";           Dif = 12; Seasons = @($s1213_1314) }
    @{ Row = 35; Team = "Rayo Vallecano";   Dif = 11; Seasons = @($s1415_1516) }
    @{ Row = 36; Team = "Real Madrid";      Dif = 15; Seasons = @($s1112_1213) }
    @{ Row = 37; Team = "Real Sociedad";    Dif = 19; Seasons = @($s1112_1213) }
    @{ Row = 38; Team = "Sevilla FC";       Dif = 24; Seasons = @($s1415_1516) }
    @{ Row = 39; Team = "Valencia";         Dif = 33; Seasons = @($s1415_1516) }
)

foreach ($entry in $teamRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Team
    $ws.Cells.Item($r, 2).Value = $entry.Dif
    $col = 3
    foreach ($season in $entry.Seasons) {
        $ws.Cells.Item($r, $col).Value = $season
        $col++
    }
}

Write-Output "applied Punto 7 Programa 1 edit"
